$p = $ppt.ActivePresentation

# --- Change 1: Slide 5 ("3a Estoria") -------------------------------------
# Fix missing space before the opening parenthesis in the "Descricao" line.
$s5 = $p.Slides.Item(5)
$sh5 = $s5.Shapes.Item(2)
$tr5 = $sh5.TextFrame.TextRange
$tr5.Paragraphs(2).Runs(2).Text = ": Adicionar fornecedor, gerar relatórios (todos, só entrada ou só saída)."

# --- Change 2: Slide 7 ("1a Estoria") --------------------------------------
# Add a new "Status: [Em andamento]." paragraph after the "Descricao" line.
$s7 = $p.Slides.Item(7)
$sh7 = $s7.Shapes.Item(2)
$tr7 = $sh7.TextFrame.TextRange

$inserted = $tr7.InsertAfter("`rStatus: [Em andamento].")
$newPara = $tr7.Paragraphs($tr7.Paragraphs().Count)
$newPara.ParagraphFormat.Alignment = 4

$newPara.Characters(1, 9).Font.Bold = -1
$newPara.Characters(22, 1).Font.Bold = -1
